$wb = $excel.ActiveWorkbook

# ---------- OpsTracker sheet ----------
$ws1 = $wb.Worksheets.Item("OpsTracker")

# Update owners for a few existing rows
$ws1.Cells.Item(12,3).Value = "Debasish"   # row12: Teaching Ops training to team
$ws1.Cells.Item(16,3).Value = "Victor"     # row16: Stamp Pad
$ws1.Cells.Item(24,3).Value = "Debasish"   # row24: Student fees collection report
$ws1.Cells.Item(27,3).Value = "Victor"     # row27: ID card for all

# Add comment to row 29 (Buy JELET Book)
$ws1.Cells.Item(29,5).Value = "Book has to given to Avishek and Subroto Sir"

# Row32: change text + item number
$ws1.Cells.Item(32,1).Value = 32
$ws1.Cells.Item(32,2).Value = "Register for CRM information"

# New rows 33 and 34
$ws1.Cells.Item(33,1).Value = 33
$ws1.Cells.Item(33,2).Value = "Register for bio data of teacher"
$ws1.Cells.Item(33,3).Value = "Victor"
$ws1.Cells.Item(33,4).Value = "Todo"

$ws1.Cells.Item(34,1).Value = 34
$ws1.Cells.Item(34,2).Value = "Organize the office different accessories"
$ws1.Cells.Item(34,3).Value = "Victor"
$ws1.Cells.Item(34,4).Value = "Todo"

# ---------- InternalAdmin sheet ----------
$ws3 = $wb.Worksheets.Item("InternalAdmin")
$ws3.Cells.Item(11,3).Value = "Done"  # Attendance Register status Todo->Done (row 11 before delete)
$ws3.Rows.Item(7).Delete()  # delete the Biscuits row

# Renumber the Item Number column (A) for rows after the deleted one
$ws3.Cells.Item(7,1).Value = 6
$ws3.Cells.Item(8,1).Value = 7
$ws3.Cells.Item(9,1).Value = 8
$ws3.Cells.Item(10,1).Value = 9
$ws3.Cells.Item(11,1).Value = 10

Write-Host "All edits applied"
